{"js": "// The document describes the 'sub-task' branch naming convention as\n// 'sub-task/<featureNummer>-<subtaskNaam>'. This edit inserts the literal\n// text \"feature-\" right before \"<featureNummer\", turning it into\n// 'sub-task/feature-<featureNummer>-<subtaskNaam>'.\n//\n// We anchor the search on the full, unambiguous phrase \"sub-task/<featureNummer>\"\n// so we don't accidentally touch the very similar 'bugfix/<featureNummer>-<bug>'\n// bullet elsewhere in the document.\nconst body = context.document.body;\n\nconst anchors = body.search(\"sub-task/<featureNummer\", { matchCase: true });\nanchors.load(\"text\");\nawait context.sync();\n\nif (anchors.items.length === 0) {\n  throw new Error(\"Could not find the 'sub-task/<featureNummer' anchor text\");\n}\n\n// Within that unique match, locate the \"<\" that starts \"<featureNummer\" and\n// insert \"feature-\" immediately before it (i.e. right after the \"/\").\nconst anchorRange = anchors.items[0];\nconst ltMatches = anchorRange.search(\"<\", { matchCase: true });\nltMatches.load(\"text\");\nawait context.sync();\n\nif (ltMatches.items.length === 0) {\n  throw new Error(\"Could not find '<' inside the anchor range\");\n}\n\nltMatches.items[0].insertText(\"feature-\", \"Before\");\nawait context.sync();\n", "ps1": "# The document describes the 'sub-task' branch naming convention as\n# 'sub-task/<featureNummer>-<subtaskNaam>'. This edit inserts the literal\n# text \"feature-\" right before \"<featureNummer\", turning it into\n# 'sub-task/feature-<featureNummer>-<subtaskNaam>'.\n#\n# We anchor the Find on the full, unambiguous phrase \"sub-task/<featureNummer\"\n# so we don't accidentally touch the very similar 'bugfix/<featureNummer>-<bug>'\n# bullet elsewhere in the document.\n\n$wdCollapseStart = 1\n\n$d = $word.ActiveDocument\n\n$anchor = $d.Content\n$found = $anchor.Find.Execute(\"sub-task/<featureNummer\")\n\nif (-not $found) {\n    throw \"Could not find the 'sub-task/<featureNummer' anchor text\"\n}\n\n# Work on a copy of the matched range so we don't disturb $anchor itself;\n# narrow it down to the \"<\" that starts \"<featureNummer\" and collapse to a\n# caret immediately before it (i.e. right after the \"/\"), then type the\n# missing \"feature-\" there.\n$target = $anchor.Duplicate\n$ltFound = $target.Find.Execute(\"<\")\n\nif (-not $ltFound) {\n    throw \"Could not find '<' inside the anchor range\"\n}\n\n$target.Collapse($wdCollapseStart)\n$target.InsertBefore(\"feature-\")\n"}
